$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.02258322285507441
$ws.Range("C2").Value = 0.0001537489499301437
$ws.Range("D2").Value = 3.082599426703578
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("G2").Value = 3.60522310558264

$ws.Range("B3").Value = 1.505614041169197
$ws.Range("C3").Value = 0.3375848360084654
$ws.Range("D3").Value = 0.7127328510149897
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 3.055818435266709
